$wb = $excel.ActiveWorkbook
$settings = $wb.Worksheets.Item("Settings")

# 1. Create the Email sheet by copying Settings (to preserve row/style fidelity), then clean it up
$settings.Copy($null, $settings)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "Email"

foreach ($h in $newSheet.Hyperlinks) {
    $h.Delete()
}
$newSheet.Range("D1:G3").Clear() | Out-Null
$newSheet.Range("G2").Clear() | Out-Null
$newSheet.Range("B3").Clear() | Out-Null

$newSheet.Range("A2").Value = "email"
$newSheet.Range("B2").Value = "diana.gradinaru.sincai@gmail.com"
$newSheet.Range("C2").Value = "Email address"

$newSheet.Range("A3").Value = "subject"
$newSheet.Range("C3").Value = "Email subject"
$newSheet.Range("B3").Value = "Assignment number 1 done"

$newSheet.Range("B4").WrapText = $true

$newSheet.Columns.Item(1).ColumnWidth = 7.42578125
$newSheet.Columns.Item(2).ColumnWidth = 80.42578125
$newSheet.Columns.Item(3).ColumnWidth = 18.85546875

$newSheet.Range("B17").Select()

# 2. Update the Settings sheet itself: remove hyperlink, update url value/description, reselect D16
foreach ($h in $settings.Hyperlinks) {
    $h.Delete()
}
$settings.Range("G2").Value = $null
$settings.Range("B2").Value = "https://www.mobile.de/ro/automobil/mazda-cx-5/vhc:car,cnt:de,pgn:1,pgs:50,ms1:16800_33_,frn:2012,frx:2018,ful:diesel!electricity,mlx:100000"
$settings.Range("C2").Value = "url for the filtered data, 50 results per page"

$settings.Select()
$settings.Range("D16").Select()
